# Update cryptos list with refreshed prices / volume(1h) figures.
# Every target cell is a plain-text cell in the source workbook (t="inlineStr").
# Setting .Value on a numeric-looking string (e.g. "1.01") would otherwise let
# Excel auto-coerce it to a Number, so each write is wrapped with a text
# NumberFormat ("@") and a ClearFormats() afterwards, which restores the
# cell to its original (default) style while keeping the value as Text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "42.492.57"
Set-TextValue "E2" "  -3.00%  "
Set-TextValue "D3" "2.228.61"
Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "111.86"
Set-TextValue "E5" "  -6.88%  "
Set-TextValue "D6" "296.97"
Set-TextValue "D7" "0.626"
Set-TextValue "E7" "  -3.13%  "
Set-TextValue "E8" "  -0.13%  "
Set-TextValue "D9" "0.608"
Set-TextValue "E9" "  -3.44%  "
Set-TextValue "D10" "44.41"
Set-TextValue "E10" "  -8.57%  "
Set-TextValue "D11" "0.0917"
Set-TextValue "E11" "  -3.63%  "
Set-TextValue "D12" "54.31"
Set-TextValue "E12" "  +0.02%  "
Set-TextValue "E13" "  -4.44%  "
Set-TextValue "E14" "  +9.85%  "
Set-TextValue "E15" "  -2.84%  "
Set-TextValue "E16" "  -3.37%  "
Set-TextValue "D17" "2.560.31"
Set-TextValue "E17" "  -2.57%  "
Set-TextValue "D18" "2.236.30"
Set-TextValue "E18" "  -2.14%  "
Set-TextValue "D19" "42.468.90"
Set-TextValue "E19" "  -3.06%  "
Set-TextValue "D20" "7.41"
Set-TextValue "E20" "  +6.46%  "
Set-TextValue "E21" "  -4.32%  "
Set-TextValue "D22" "72.76"
Set-TextValue "E22" "  +0.50%  "
Set-TextValue "D23" "3.50"
Set-TextValue "E23" "  +21.39%  "
Set-TextValue "E24" "  -2.11%  "
Set-TextValue "D25" "229.64"
Set-TextValue "E25" "  -2.96%  "
Set-TextValue "D26" "9.25"
Set-TextValue "E26" "  -4.09%  "
Set-TextValue "D27" "11.74"
Set-TextValue "E27" "  -2.56%  "
Set-TextValue "E28" "  -1.76%  "
Set-TextValue "D29" "2.24"
Set-TextValue "E29" "  -0.67%  "
Set-TextValue "D30" "38.49"
Set-TextValue "E30" "  -9.81%  "
Set-TextValue "D31" "3.25"
Set-TextValue "E31" "  -3.99%  "
Set-TextValue "D32" "174.31"
Set-TextValue "E32" "  +0.99%  "
Set-TextValue "D33" "21.09"
Set-TextValue "E33" "  -2.79%  "
Set-TextValue "D34" "0.0895"
Set-TextValue "E34" "  -4.09%  "
Set-TextValue "D35" "5.22"
Set-TextValue "E35" "  +12.87%  "
Set-TextValue "D36" "5.69"
Set-TextValue "E36" "  -1.88%  "
Set-TextValue "D37" "4.38"
Set-TextValue "E37" "  +2.21%  "
Set-TextValue "E38" "  -3.73%  "
Set-TextValue "E39" "  -1.51%  "
Set-TextValue "E40" "  -1.78%  "
Set-TextValue "E41" "  -5.05%  "
Set-TextValue "D42" "72.14"
Set-TextValue "E42" "  -3.09%  "
Set-TextValue "E43" "  -0.84%  "
Set-TextValue "E44" "  -7.32%  "
Set-TextValue "E45" "  +0.14%  "
Set-TextValue "D46" "1.33"
Set-TextValue "E46" "  -4.91%  "
Set-TextValue "D47" "5.48"
Set-TextValue "E47" "  -6.46%  "
Set-TextValue "D48" "1.32"
Set-TextValue "E48" "  +4.11%  "
Set-TextValue "D49" "103.29"
Set-TextValue "E49" "  +0.71%  "
Set-TextValue "B50" "Stacks"
Set-TextValue "C50" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D50" "1.66"
Set-TextValue "E50" "  +7.34%  "
Set-TextValue "B51" "FraxShare"
Set-TextValue "C51" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D51" "8.49"
Set-TextValue "E51" "  -0.96%  "
